$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "citizenship-with-scenario-outline;citizenship-create-and-delete-with-scenario-outline;;2"
$ws.Range("B61").Value = "passed"
$ws.Range("C61").Value = "chrome"
$ws.Range("D61").Value = "21.09.2021"

$ws.Range("A62").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B62").Value = "passed"
$ws.Range("C62").Value = "chrome"
$ws.Range("D62").Value = "26.09.2021"

$ws.Range("A63").Value = "country-multi-scenario;create-a-country"
$ws.Range("B63").Value = "failed"
$ws.Range("C63").Value = "chrome"
$ws.Range("D63").Value = "26.09.2021"

$ws.Range("A64").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B64").Value = "failed"
$ws.Range("C64").Value = "chrome"
$ws.Range("D64").Value = "26.09.2021"

$ws.Range("A65").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B65").Value = "passed"
$ws.Range("C65").Value = "chrome"
$ws.Range("D65").Value = "26.09.2021"

$ws.Range("A66").Value = "country-functionality;create-country"
$ws.Range("B66").Value = "passed"
$ws.Range("C66").Value = "chrome"
$ws.Range("D66").Value = "26.09.2021"

$ws.Range("A67").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B67").Value = "passed"
$ws.Range("C67").Value = "chrome"
$ws.Range("D67").Value = "26.09.2021"

$ws.Range("A68").Value = "country-multi-scenario;create-a-country"
$ws.Range("B68").Value = "passed"
$ws.Range("C68").Value = "chrome"
$ws.Range("D68").Value = "26.09.2021"

$ws.Range("A69").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B69").Value = "passed"
$ws.Range("C69").Value = "chrome"
$ws.Range("D69").Value = "26.09.2021"

$ws.Range("A70").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B70").Value = "passed"
$ws.Range("C70").Value = "chrome"
$ws.Range("D70").Value = "26.09.2021"

$ws.Range("A71").Value = "country-multi-scenario;create-a-country"
$ws.Range("B71").Value = "passed"
$ws.Range("C71").Value = "chrome"
$ws.Range("D71").Value = "26.09.2021"

$ws.Range("A72").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B72").Value = "failed"
$ws.Range("C72").Value = "chrome"
$ws.Range("D72").Value = "26.09.2021"

$ws.Range("A73").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B73").Value = "passed"
$ws.Range("C73").Value = "chrome"
$ws.Range("D73").Value = "26.09.2021"

$ws.Range("A74").Value = "country-multi-scenario;create-a-country"
$ws.Range("B74").Value = "failed"
$ws.Range("C74").Value = "chrome"
$ws.Range("D74").Value = "26.09.2021"

$ws.Range("A75").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B75").Value = "passed"
$ws.Range("C75").Value = "chrome"
$ws.Range("D75").Value = "26.09.2021"

$ws.Range("A76").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B76").Value = "passed"
$ws.Range("C76").Value = "chrome"
$ws.Range("D76").Value = "26.09.2021"

$ws.Range("A77").Value = "country-multi-scenario;create-a-country"
$ws.Range("B77").Value = "passed"
$ws.Range("C77").Value = "chrome"
$ws.Range("D77").Value = "26.09.2021"

$ws.Range("A78").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B78").Value = "passed"
$ws.Range("C78").Value = "chrome"
$ws.Range("D78").Value = "26.09.2021"

$ws.Range("A79").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B79").Value = "passed"
$ws.Range("C79").Value = "chrome"
$ws.Range("D79").Value = "26.09.2021"

$ws.Range("A80").Value = "login-functionality;login-with-valid-username-and-password"
$ws.Range("B80").Value = "passed"
$ws.Range("C80").Value = "chrome"
$ws.Range("D80").Value = "26.09.2021"

$ws.Range("A81").Value = "country-multi-scenario;create-a-country"
$ws.Range("B81").Value = "passed"
$ws.Range("C81").Value = "chrome"
$ws.Range("D81").Value = "26.09.2021"

$ws.Range("A82").Value = "country-multi-scenario;create-a-country-parameter-data"
$ws.Range("B82").Value = "passed"
$ws.Range("C82").Value = "chrome"
$ws.Range("D82").Value = "26.09.2021"

